$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..25, columns A..F (1..6), are being reordered (a pure row
# permutation - the header row 1 and the totals row 26 are untouched).
# newRow[r] = oldRow[mapping[r]]
$mapping = @{
    2  = 15
    3  = 5
    4  = 7
    5  = 13
    6  = 14
    7  = 6
    8  = 3
    9  = 4
    10 = 10
    11 = 9
    12 = 12
    13 = 11
    14 = 2
    15 = 8
    16 = 19
    17 = 16
    18 = 20
    19 = 17
    20 = 18
    21 = 21
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

# Snapshot every source row's values before any writes happen, so that
# writes to earlier rows don't clobber data still needed for later rows.
$snapshot = @{}
foreach ($r in $mapping.Values) {
    if (-not $snapshot.ContainsKey($r)) {
        $rowVals = @()
        for ($c = 1; $c -le 6; $c++) {
            $rowVals += , ($ws.Cells.Item($r, $c).Value())
        }
        $snapshot[$r] = $rowVals
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
